$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "size"
$ws.Range("B1").Value = "color"
$ws.Range("C1").Value = ""
